# Updates cryptos list values (price + 1h volume change) per upstream source refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text looks like a plain number need an explicit text format first,
# otherwise Excel COM auto-coerces the assigned string into a numeric cell value.

$ws.Range('D2').Value = '30.320.90'
$ws.Range('E2').Value = '  -3.70%  '

$ws.Range('D3').Value = '1.930.58'
$ws.Range('E3').Value = '  -3.91%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9996'
$ws.Range('E4').Value = '  +0.00%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '248.65'
$ws.Range('E5').Value = '  -4.32%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.7262'
$ws.Range('E6').Value = '  -5.51%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.9997'
$ws.Range('E7').Value = '  +0.06%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3263'
$ws.Range('E8').Value = '  -9.09%  '

$ws.Range('E9').Value = '  -4.28%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.06805'
$ws.Range('E10').Value = '  -3.81%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.8022'
$ws.Range('E11').Value = '  -4.76%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.08057'
$ws.Range('E12').Value = '  -0.68%  '

$ws.Range('D13').Value = '1.930.86'
$ws.Range('E13').Value = '  -3.87%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.409'
$ws.Range('E14').Value = '  -4.02%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '94.67'
$ws.Range('E15').Value = '  -6.66%  '

$ws.Range('E16').Value = '  -1.51%  '

$ws.Range('D17').Value = '30.305.77'
$ws.Range('E17').Value = '  -3.78%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '253.88'
$ws.Range('E18').Value = '  -7.73%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.000008005'
$ws.Range('E19').Value = '  +0.07%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '5.823'
$ws.Range('E20').Value = '  -2.23%  '

$ws.Range('D21').Value = '2.184.00'
$ws.Range('E21').Value = '  -3.72%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.9996'
$ws.Range('E22').Value = '  +0.04%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.9992'
$ws.Range('E23').Value = '  -0.11%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '6.861'
$ws.Range('E24').Value = '  -5.07%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '9.676'
$ws.Range('E25').Value = '  -4.74%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '159.03'
$ws.Range('E26').Value = '  -3.16%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.374'
$ws.Range('E27').Value = '  -0.67%  '

$ws.Range('B28').Value = 'Stellar'
$ws.Range('C28').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.1342'
$ws.Range('E28').Value = '  -7.91%  '

$ws.Range('B29').Value = 'EthereumClassic'
$ws.Range('C29').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '19.05'
$ws.Range('E29').Value = '  -5.58%  '

$ws.Range('E30').Value = '  -4.54%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.339'
$ws.Range('E31').Value = '  -1.38%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.396'
$ws.Range('E32').Value = '  -5.29%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.185'
$ws.Range('E33').Value = '  -4.98%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.05060'
$ws.Range('E34').Value = '  -3.02%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.214'
$ws.Range('E35').Value = '  -2.06%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.7382'
$ws.Range('E36').Value = '  -3.30%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.754'
$ws.Range('E37').Value = '  -1.67%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01967'
$ws.Range('E38').Value = '  -2.68%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.816'
$ws.Range('E39').Value = '  -4.52%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '6.587'
$ws.Range('E40').Value = '  -2.04%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '79.00'
$ws.Range('E41').Value = '  -1.72%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.4454'
$ws.Range('E42').Value = '  -6.42%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.987'
$ws.Range('E43').Value = '  -9.34%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.000'
$ws.Range('E44').Value = '  +0.06%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.8337'
$ws.Range('E45').Value = '  -3.65%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '101.86'
$ws.Range('E46').Value = '  -2.78%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '9.727'
$ws.Range('E47').Value = '  -2.28%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '7.273'
$ws.Range('E48').Value = '  -5.39%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '36.30'
$ws.Range('E49').Value = '  -2.27%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.05933'
$ws.Range('E50').Value = '  -0.78%  '

$ws.Range('B51').Value = 'Decentraland'
$ws.Range('C51').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.4068'
$ws.Range('E51').Value = '  -7.15%  '
